$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row: rename "Phone_Number" -> "phone" (others stay the same value)
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "phone"

# ---------------------------------------------------------------------------
# Existing data row (row 2): prefix the phone number with the country code
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 919510038048

# ---------------------------------------------------------------------------
# New row 3: Mohit Aswani
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Mohit Aswani "
$ws.Range("B3").Value = 919328027733
$ws.Range("B3").Font.Name = "Segoe UI"
$ws.Range("B3").Font.Color = 2956306
$ws.Range("C3").Value = "mohit.creerinfotech@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:mohit.creerinfotech@gmail.com") | Out-Null
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D3").Value = "Ahmedabad "
$ws.Rows.Item(3).RowHeight = 16.5

# ---------------------------------------------------------------------------
# New row 4: Shreeji Nandola
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Shreeji Nandola"
$ws.Range("B4").Value = 919106284482
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C4").Value = "mohit.creerinfotech@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:mohit.creerinfotech@gmail.com") | Out-Null
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D4").Value = "Ahmedabad "
$ws.Rows.Item(4).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Page setup (printer friendly portrait / A4-ish "9" paper size)
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Move the selection back to the top of the sheet
# ---------------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
